$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 2 ("H 72") - data shifts up, shrinking the used range from
# A1:F63 to A1:F62.
$ws.Rows.Item(2).Delete()
